$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.910.11"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.917.97"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "324.66"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.4567"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").Value = "0.3798"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "0.07741"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").Value = "0.9751"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").Value = "22.26"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "1.933.58"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "5.689"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "6.961"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "0.06992"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "84.31"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "0.000009455"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "16.63"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D21").Value = "28.925.78"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "5.331"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").Value = "2.135.13"
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").Value = "157.60"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "18.95"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").Value = "5.601"
$ws.Range("D29").Value = "117.57"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").Value = "0.09298"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "0.8649"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "5.089"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "1.238"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "3.014"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "0.02037"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").Value = "3.065"
$ws.Range("E40").Value = "  +11.67%  "
$ws.Range("D41").Value = "7.460"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "0.5480"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").Value = "0.1751"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "9.309"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "0.000002767"
$ws.Range("E45").Value = "  +16.72%  "
$ws.Range("D46").Value = "2.156"
$ws.Range("E46").Value = "  +3.38%  "
$ws.Range("D47").Value = "0.5148"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "0.06936"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "11.09"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").Value = "110.48"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("E51").Value = "  -0.78%  "
